# Auto-generated edit script applying numeric corrections to Halicarnassus_Profits sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33 (item id 5512), hunk 0
$ws.Range("H33").Value = 243.16667
$ws.Range("J33").Value = 75
$ws.Range("L33").Value = 75
$ws.Range("N33").Value = -533
# Row 45 (item id 4585), hunk 1
$ws.Range("H45").Value = 1176.2222
$ws.Range("I45").Value = 99
$ws.Range("J45").Value = 1484
$ws.Range("K45").Value = 297
$ws.Range("L45").Value = 4452
$ws.Range("M45").Value = -105
$ws.Range("N45").Value = -4836
# Row 86 (item id 12603), hunk 2
$ws.Range("H86").Value = 2980.1
$ws.Range("I86").Value = 2543.1428
$ws.Range("K86").Value = 2543.1428
$ws.Range("M86").Value = -1420.1428
# Row 89 (item id 12603), hunk 3
$ws.Range("H89").Value = 2980.1
$ws.Range("I89").Value = 2543.1428
$ws.Range("K89").Value = 12715.714
$ws.Range("M89").Value = -7099.714
# Row 111 (item id 27768), hunk 4
$ws.Range("H111").Value = 1294.3334
$ws.Range("J111").Value = 1444
$ws.Range("L111").Value = 4332
$ws.Range("N111").Value = -10466
# Row 132 (item id 44049), hunk 5
$ws.Range("H132").Value = 18591.584
$ws.Range("I132").Value = 18591.584
$ws.Range("K132").Value = 55774.75199999999
$ws.Range("M132").Value = -53244.75199999999
# Row 135 (item id 44047), hunk 6
$ws.Range("H135").Value = 541
$ws.Range("I135").Value = 541
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 4869
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -2334
$ws.Range("N135").ClearContents()
# Row 141 (item id 44161), hunk 7
$ws.Range("H141").Value = 4153
$ws.Range("I141").Value = 4448.5
$ws.Range("K141").Value = 13345.5
$ws.Range("M141").Value = -8165.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (item id 27713), hunk 8
$ws.Range("H2").Value = 897.53845
$ws.Range("I2").Value = 897.53845
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 897.53845
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -784.53845
$ws.Range("N2").ClearContents()
# Row 45 (item id 27714), hunk 9
$ws.Range("H45").Value = 2667.9092
$ws.Range("I45").Value = 2279.8572
$ws.Range("K45").Value = 2279.8572
$ws.Range("M45").Value = -1902.8572
# Row 94 (item id 18055), hunk 10
$ws.Range("H94").Value = 4330
$ws.Range("J94").Value = 4330
$ws.Range("L94").Value = 4330
$ws.Range("N94").Value = -6132
# Row 97 (item id 19941), hunk 11
$ws.Range("H97").Value = 867.64703
$ws.Range("I97").Value = 836.6667
$ws.Range("J97").Value = 1100
$ws.Range("K97").Value = 836.6667
$ws.Range("L97").Value = 1100
$ws.Range("M97").Value = -340.6667
$ws.Range("N97").Value = -2092
# Row 102 (item id 19945), hunk 12
$ws.Range("H102").Value = 4432.6665
$ws.Range("I102").Value = 3040.8333
$ws.Range("K102").Value = 3040.8333
$ws.Range("M102").Value = -1418.8333
# Row 109 (item id 25646), hunk 13
$ws.Range("H109").Value = 150001
$ws.Range("J109").Value = 150001
$ws.Range("L109").Value = 150001
$ws.Range("N109").Value = -152775
# Row 110 (item id 27708), hunk 14
$ws.Range("H110").Value = 803.8570999999999
$ws.Range("J110").Value = 727
$ws.Range("L110").Value = 727
$ws.Range("N110").Value = -4817
# Row 116 (item id 27713), hunk 15
$ws.Range("H116").Value = 897.53845
$ws.Range("I116").Value = 897.53845
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 897.53845
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1396.46155
$ws.Range("N116").ClearContents()
# Row 122 (item id 36168), hunk 16
$ws.Range("H122").Value = 2898.8
$ws.Range("I122").Value = 2498.5
$ws.Range("K122").Value = 7495.5
$ws.Range("M122").Value = -5045.5
# Row 132 (item id 43997), hunk 17
$ws.Range("H132").Value = 4385.68
$ws.Range("I132").Value = 3482.1
$ws.Range("K132").Value = 10446.3
$ws.Range("M132").Value = -7916.299999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (item id 27713), hunk 18
$ws.Range("H3").Value = 897.53845
$ws.Range("I3").Value = 897.53845
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 897.53845
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -783.53845
$ws.Range("N3").ClearContents()
# Row 7 (item id 1602), hunk 19
$ws.Range("H7").Value = 10714457
$ws.Range("I7").Value = 11000195
$ws.Range("J7").Value = 10000112
$ws.Range("K7").Value = 11000195
$ws.Range("L7").Value = 10000112
$ws.Range("M7").Value = -11000082
$ws.Range("N7").Value = -10000338
# Row 20 (item id 14149), hunk 20
$ws.Range("H20").Value = 972
$ws.Range("I20").Value = 972
$ws.Range("K20").Value = 972
$ws.Range("M20").Value = -725
# Row 86 (item id 12526), hunk 21
$ws.Range("H86").Value = 7278.3335
$ws.Range("I86").Value = 5472.25
$ws.Range("J86").Value = 9342.429
$ws.Range("K86").Value = 5472.25
$ws.Range("L86").Value = 9342.429
$ws.Range("M86").Value = -4349.25
$ws.Range("N86").Value = -11588.429
# Row 89 (item id 12526), hunk 22
$ws.Range("H89").Value = 7278.3335
$ws.Range("I89").Value = 5472.25
$ws.Range("J89").Value = 9342.429
$ws.Range("K89").Value = 27361.25
$ws.Range("L89").Value = 46712.145
$ws.Range("M89").Value = -21745.25
$ws.Range("N89").Value = -57944.145
# Row 94 (item id 19939), hunk 23
$ws.Range("H94").Value = 1831.5555
$ws.Range("I94").Value = 1848
$ws.Range("K94").Value = 1848
$ws.Range("M94").Value = -1397
# Row 99 (item id 19943), hunk 24
$ws.Range("H99").Value = 1829.8334
$ws.Range("I99").Value = 1777
$ws.Range("K99").Value = 1777
$ws.Range("M99").Value = -279

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 99 (item id 36198), hunk 25
$ws.Range("H99").Value = 2400
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
# Row 126 (item id 36198), hunk 26
$ws.Range("H126").Value = 2400
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 107 (item id 27838), hunk 27
$ws.Range("H107").Value = 369.25
$ws.Range("I107").Value = 367.66666
$ws.Range("J107").Value = 374
$ws.Range("K107").Value = 1102.99998
$ws.Range("L107").Value = 1122
$ws.Range("M107").Value = 817.0000199999999
$ws.Range("N107").Value = -4962
# Row 122 (item id 36078), hunk 28
$ws.Range("H122").Value = 983
$ws.Range("I122").Value = 974.5
$ws.Range("K122").Value = 8770.5
$ws.Range("M122").Value = -6320.5
# Row 131 (item id 36060), hunk 29
$ws.Range("H131").Value = 1561.6
$ws.Range("J131").Value = 2420
$ws.Range("L131").Value = 7260
$ws.Range("N131").Value = -17340

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2 (item id 5062), hunk 30
$ws.Range("H2").Value = 194.16667
$ws.Range("I2").Value = 46.545456
$ws.Range("J2").Value = 426.14285
$ws.Range("K2").Value = 46.545456
$ws.Range("L2").Value = 426.14285
$ws.Range("M2").Value = 66.454544
$ws.Range("N2").Value = -652.14285
# Row 35 (item id 4317), hunk 31
$ws.Range("H35").Value = 5000
$ws.Range("I35").Value = 5000
$ws.Range("K35").Value = 5000
$ws.Range("M35").Value = -4702
# Row 51 (item id 27222), hunk 32
$ws.Range("H51").Value = 92000
$ws.Range("J51").Value = 92000
$ws.Range("L51").Value = 92000
$ws.Range("N51").Value = -93018
# Row 70 (item id 14146), hunk 33
$ws.Range("H70").Value = 11007
$ws.Range("I70").Value = 11007
$ws.Range("K70").Value = 11007
$ws.Range("M70").Value = -10737
# Row 73 (item id 14146), hunk 34
$ws.Range("H73").Value = 11007
$ws.Range("I73").Value = 11007
$ws.Range("K73").Value = 11007
$ws.Range("M73").Value = -10071
# Row 80 (item id 12521), hunk 35
$ws.Range("H80").Value = 2254.818
$ws.Range("I80").Value = 1841.75
$ws.Range("K80").Value = 1841.75
$ws.Range("M80").Value = -843.75
# Row 83 (item id 12521), hunk 36
$ws.Range("H83").Value = 2254.818
$ws.Range("I83").Value = 1841.75
$ws.Range("K83").Value = 9208.75
$ws.Range("M83").Value = -4216.75
# Row 102 (item id 36169), hunk 37
$ws.Range("H102").Value = 1866.1428
$ws.Range("I102").Value = 1728.7273
$ws.Range("J102").Value = 2370
$ws.Range("K102").Value = 1728.7273
$ws.Range("L102").Value = 2370
$ws.Range("M102").Value = -106.7273
$ws.Range("N102").Value = -5614
# Row 107 (item id 27802), hunk 38
$ws.Range("H107").Value = 1232.0588
$ws.Range("I107").Value = 1457.375
$ws.Range("J107").Value = 1031.7778
$ws.Range("K107").Value = 1457.375
$ws.Range("L107").Value = 1031.7778
$ws.Range("M107").Value = 462.625
$ws.Range("N107").Value = -4871.7778
# Row 129 (item id 35367), hunk 39
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("N129").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22 (item id 5277), hunk 40
$ws.Range("H22").Value = 5000
$ws.Range("J22").Value = 5000
$ws.Range("L22").Value = 5000
$ws.Range("N22").Value = -5590
# Row 27 (item id 5277), hunk 41
$ws.Range("H27").Value = 5000
$ws.Range("J27").Value = 5000
$ws.Range("L27").Value = 5000
$ws.Range("N27").Value = -5214
# Row 40 (item id 36248), hunk 42
$ws.Range("H40").Value = 8000
$ws.Range("I40").Value = 8000
$ws.Range("K40").Value = 8000
$ws.Range("M40").Value = -7864
# Row 93 (item id 19993), hunk 43
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
# Row 122 (item id 36247), hunk 44
$ws.Range("H122").Value = 3500
$ws.Range("I122").Value = 3500
$ws.Range("K122").Value = 10500
$ws.Range("M122").Value = -8050

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 96 (item id 19977), hunk 45
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
# Row 112 (item id 25836), hunk 46
$ws.Range("H112").Value = 36965.668
$ws.Range("J112").Value = 36965.668
$ws.Range("L112").Value = 36965.668
$ws.Range("N112").Value = -39919.668

